$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the Agra branch address (remove the "xx" typo before the trailing period)
$ws.Range("B2").Value = "Mercury Travels Ltd.`nC/o Hotel Clarks Shiraz, 54, Taj Road, Agra 282 001."

# Ahmedabad branch address stays the same text
$ws.Range("B3").Value = "Mercury Travels Ltd.`n103 B Abhijeet 1, 1st Floor, Mithakali Six Roads, Ahmedabad 380 006."

# Introduce a typo ("xxx") into the Bangalore branch address
$ws.Range("B4").Value = "Mercury Travels Ltd.`n125, Infantry Road, Bangalore 560 001xxx."

# Move the active selection to B4
$ws.Range("B4").Select()
